$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values from 300 to 3000 for existing rows 1-72
for ($r = 1; $r -le 72; $r++) {
    $ws.Cells.Item($r, 1).Value = 3000
}

# Add two new rows (73, 74) with A = 3000, matching the style (s="4") used
# by the rest of column A (copy formats from A1, then set the value)
$ws.Range("A1").Copy()
$ws.Range("A73:A74").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(73, 1).Value = 3000
$ws.Cells.Item(74, 1).Value = 3000

# Add formatted (empty) cells S63:S64, matching style of existing N12 (s="3")
$ws.Range("N12").Copy()
$ws.Range("S63:S64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view: scroll so A37 is the top-left cell, and select S63:S64
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("S63:S64").Select()
